$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds header "K" (formerly based on "Strike#"); regenerate values.
$ws.Range("G2").Value = 3
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 3
$ws.Range("G6").Value = 0
$ws.Range("G7").Value = 2
$ws.Range("G8").Value = 2
$ws.Range("G10").Value = 2
